# Simulated Wild Card round and logged it
# Appends new play-by-play data points to the running season logs (YDS and ST
# sheets hold space-separated sample lists in shared strings) and updates the
# season-to-date aggregate totals on the OFF, DEF, ST, TURNS and PEN sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# YDS sheet: append new per-play yardage samples to the four running lists
# ---------------------------------------------------------------------------
$wsYDS = $wb.Worksheets.Item("YDS")

$wsYDS.Range("B2").Value2 = $wsYDS.Range("B2").Value2 + " 8 1 2 3 0 10 -2 6 5 2 8 23 1 13 8"
$wsYDS.Range("B3").Value2 = $wsYDS.Range("B3").Value2 + " 11 33 6 9 11 4 7 10 3 18 4 10 8 21 9 12"
$wsYDS.Range("C2").Value2 = $wsYDS.Range("C2").Value2 + " 3 -1 1 2 1 4 2 4 1 12 8 2 0 26 4 1 0 2 4 5 1 6 0 6 3 -1 6"
$wsYDS.Range("C3").Value2 = $wsYDS.Range("C3").Value2 + " 17 6 12 16 8 15 7 2 9 1 11 6 18 29 7 10 14 12 8 1 6 5 3"

# ---------------------------------------------------------------------------
# OFF sheet: update Home (row 2) / Road (row 3) season totals
# ---------------------------------------------------------------------------
$wsOFF = $wb.Worksheets.Item("OFF")

$wsOFF.Range("C2").Value2 = 208
$wsOFF.Range("E2").Value2 = 14
$wsOFF.Range("F2").Value2 = 86
$wsOFF.Range("G2").Value2 = 63
$wsOFF.Range("J2").Value2 = 38
$wsOFF.Range("N2").Value2 = 17
$wsOFF.Range("O2").Value2 = 28
$wsOFF.Range("P2").Value2 = 17

$wsOFF.Range("B3").Value2 = 12
$wsOFF.Range("C3").Value2 = 171
$wsOFF.Range("E3").Value2 = 39
$wsOFF.Range("F3").Value2 = 80
$wsOFF.Range("G3").Value2 = 22
$wsOFF.Range("H3").Value2 = 29
$wsOFF.Range("I3").Value2 = 59
$wsOFF.Range("J3").Value2 = 49
$wsOFF.Range("L3").Value2 = 232
$wsOFF.Range("M3").Value2 = 149
$wsOFF.Range("Q3").Value2 = 505

# ---------------------------------------------------------------------------
# DEF sheet: update Home (row 2) / Road (row 3) season totals
# ---------------------------------------------------------------------------
$wsDEF = $wb.Worksheets.Item("DEF")

$wsDEF.Range("B2").Value2 = 7
$wsDEF.Range("C2").Value2 = 196
$wsDEF.Range("F2").Value2 = 66
$wsDEF.Range("G2").Value2 = 57
$wsDEF.Range("I2").Value2 = 8
$wsDEF.Range("J2").Value2 = 28
$wsDEF.Range("N2").Value2 = 14
$wsDEF.Range("O2").Value2 = 21
$wsDEF.Range("P2").Value2 = 12

$wsDEF.Range("B3").Value2 = 12
$wsDEF.Range("C3").Value2 = 205
$wsDEF.Range("D3").Value2 = 4
$wsDEF.Range("E3").Value2 = 26
$wsDEF.Range("F3").Value2 = 123
$wsDEF.Range("G3").Value2 = 37
$wsDEF.Range("H3").Value2 = 25
$wsDEF.Range("I3").Value2 = 70
$wsDEF.Range("J3").Value2 = 61
$wsDEF.Range("L3").Value2 = 264
$wsDEF.Range("M3").Value2 = 173
$wsDEF.Range("Q3").Value2 = 469

# ---------------------------------------------------------------------------
# ST sheet: update counts (row 2) / touchbacks (row 3), plus append new
# distance/return samples to the running lists
# ---------------------------------------------------------------------------
$wsST = $wb.Worksheets.Item("ST")

$wsST.Range("B2").Value2 = 93
$wsST.Range("D2").Value2 = 58
$wsST.Range("H2").Value2 = 6
$wsST.Range("I2").Value2 = 2
$wsST.Range("L2").Value2 = 36
$wsST.Range("M2").Value2 = 27

$wsST.Range("B3").Value2 = 47

$wsST.Range("B6").Value2 = $wsST.Range("B6").Value2 + " 40 18 41"
$wsST.Range("D3").Value2 = $wsST.Range("D3").Value2 + " 45 33"
$wsST.Range("D4").Value2 = $wsST.Range("D4").Value2 + " 0 0"
$wsST.Range("D5").Value2 = $wsST.Range("D5").Value2 + " 0 0 7"

# ---------------------------------------------------------------------------
# TURNS sheet: update Road (row 3) season totals
# ---------------------------------------------------------------------------
$wsTURNS = $wb.Worksheets.Item("TURNS")

$wsTURNS.Range("B3").Value2 = 1
$wsTURNS.Range("D3").Value2 = 9
$wsTURNS.Range("E3").Value2 = 9

# ---------------------------------------------------------------------------
# PEN sheet: update penalty counts
# ---------------------------------------------------------------------------
$wsPEN = $wb.Worksheets.Item("PEN")

$wsPEN.Range("B3").Value2 = 23
$wsPEN.Range("D4").Value2 = 8
